# Auto-generated edit script: apply updated transition-matrix probabilities
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1875
$ws.Range("C2").Value = 0.5982142857142857
$ws.Range("J2").Value = 0.01339285714285714
$ws.Range("P2").Value = 0.125
$ws.Range("S2").Value = 0.07589285714285714
$ws.Range("J3").Value = 0.01492537313432836
$ws.Range("P3").Value = 0.7910447761194029
$ws.Range("S3").Value = 0.1940298507462687
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.59375
$ws.Range("S4").Value = 0.34375
$ws.Range("B6").Value = 0.0391304347826087
$ws.Range("D6").Value = 0.02173913043478261
$ws.Range("F6").Value = 0.08695652173913043
$ws.Range("J6").Value = 0.2434782608695652
$ws.Range("O6").Value = 0.01739130434782609
$ws.Range("Q6").Value = 0.1782608695652174
$ws.Range("R6").Value = 0.0782608695652174
$ws.Range("S6").Value = 0.3347826086956522
$ws.Range("B7").Value = 0.09045226130653267
$ws.Range("E7").Value = 0.005025125628140704
$ws.Range("F7").Value = 0.08542713567839195
$ws.Range("J7").Value = 0.1306532663316583
$ws.Range("O7").Value = 0.05025125628140704
$ws.Range("Q7").Value = 0.1256281407035176
$ws.Range("R7").Value = 0.1256281407035176
$ws.Range("S7").Value = 0.3869346733668342
$ws.Range("B8").Value = 0.06775067750677506
$ws.Range("D8").Value = 0.02168021680216802
$ws.Range("E8").Value = 0.002710027100271003
$ws.Range("F8").Value = 0.06775067750677506
$ws.Range("J8").Value = 0.1327913279132791
$ws.Range("O8").Value = 0.02168021680216802
$ws.Range("Q8").Value = 0.1490514905149052
$ws.Range("R8").Value = 0.1409214092140921
$ws.Range("S8").Value = 0.3956639566395664
$ws.Range("B9").Value = 0.08256880733944955
$ws.Range("D9").Value = 0.009174311926605505
$ws.Range("F9").Value = 0.08256880733944955
$ws.Range("J9").Value = 0.1284403669724771
$ws.Range("O9").Value = 0.02752293577981652
$ws.Range("Q9").Value = 0.1467889908256881
$ws.Range("R9").Value = 0.1422018348623853
$ws.Range("S9").Value = 0.3807339449541284
$ws.Range("B10").Value = 0.08562197092084006
$ws.Range("D10").Value = 0.01453957996768982
$ws.Range("E10").Value = 0.0008077544426494346
$ws.Range("F10").Value = 0.07673667205169628
$ws.Range("J10").Value = 0.1082390953150242
$ws.Range("O10").Value = 0.02180936995153474
$ws.Range("Q10").Value = 0.2148626817447496
$ws.Range("R10").Value = 0.1138933764135703
$ws.Range("S10").Value = 0.3634894991922455
$ws.Range("G11").Value = 0.1326860841423948
$ws.Range("J11").Value = 0.06796116504854369
$ws.Range("K11").Value = 0.1941747572815534
$ws.Range("L11").Value = 0.5922330097087378
$ws.Range("S11").Value = 0.01294498381877023
$ws.Range("G12").Value = 0.734375
$ws.Range("J12").Value = 0.1979166666666667
$ws.Range("K12").Value = 0.005208333333333333
$ws.Range("L12").Value = 0.03645833333333334
$ws.Range("S12").Value = 0.02604166666666667
$ws.Range("G13").Value = 0.6764705882352942
$ws.Range("J13").Value = 0.2352941176470588
$ws.Range("S13").Value = 0.08823529411764706
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.03930131004366812
$ws.Range("H15").Value = 0.1528384279475982
$ws.Range("I15").Value = 0.05676855895196507
$ws.Range("J15").Value = 0.2882096069868996
$ws.Range("K15").Value = 0.06986899563318777
$ws.Range("M15").Value = 0.004366812227074236
$ws.Range("O15").Value = 0.07860262008733625
$ws.Range("S15").Value = 0.3100436681222707
$ws.Range("F16").Value = 0.01333333333333333
$ws.Range("H16").Value = 0.16
$ws.Range("I16").Value = 0.04666666666666667
$ws.Range("J16").Value = 0.4333333333333333
$ws.Range("K16").Value = 0.1066666666666667
$ws.Range("M16").Value = 0.006666666666666667
$ws.Range("O16").Value = 0.02666666666666667
$ws.Range("S16").Value = 0.2066666666666667
$ws.Range("F17").Value = 0.009615384615384616
$ws.Range("H17").Value = 0.2163461538461539
$ws.Range("I17").Value = 0.09615384615384616
$ws.Range("J17").Value = 0.4110576923076923
$ws.Range("K17").Value = 0.09615384615384616
$ws.Range("M17").Value = 0.009615384615384616
$ws.Range("N17").Value = 0.002403846153846154
$ws.Range("O17").Value = 0.04086538461538462
$ws.Range("S17").Value = 0.1177884615384615
$ws.Range("F18").Value = 0.01136363636363636
$ws.Range("H18").Value = 0.1590909090909091
$ws.Range("I18").Value = 0.1136363636363636
$ws.Range("J18").Value = 0.4204545454545455
$ws.Range("K18").Value = 0.1136363636363636
$ws.Range("M18").Value = 0.01893939393939394
$ws.Range("O18").Value = 0.05681818181818182
$ws.Range("S18").Value = 0.1060606060606061
$ws.Range("F19").Value = 0.0143824027072758
$ws.Range("H19").Value = 0.1548223350253807
$ws.Range("I19").Value = 0.1099830795262267
$ws.Range("J19").Value = 0.3967851099830795
$ws.Range("K19").Value = 0.1218274111675127
$ws.Range("M19").Value = 0.02030456852791878
$ws.Range("O19").Value = 0.077834179357022
$ws.Range("S19").Value = 0.1040609137055838
